$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = '[{"factor": ["serum", "thioredoxin concentration"], "outcome": ["in-hospital major adverse event"]}, {"factor": ["glasgow coma scale score", "admission"], "outcome": ["in-hospital major adverse event"]}]'
$ws.Range("G4").Value = '[{"factor": ["macrophage migration inhibitory factor", "mif"], "outcome": ["trauma", "severity", "inflammation", "clinical outcome"]}]'
$ws.Range("G5").Value = '[{"factor": ["glasgow coma scale score"], "outcome": ["gos", "glasgow outcome scale"]}, {"factor": ["artificial ventilation"], "outcome": ["mortality rate"]}, {"factor": ["monitoring", "intracranial"], "outcome": ["good", "outcome", "long-term"]}]'
$ws.Range("G6").Value = '[{"factor": ["calculator", "impact", "prognosis"], "outcome": ["elderly", "tbi", "severe", "outcome", "patient"]}, {"factor": ["expect risk", "fatal", "outcome"], "outcome": ["conservative treatment", "risk"]}, {"factor": ["outcome", "unfavorable", "predict risk"], "outcome": ["actual", "risk", "outcome", "predict", "unfavorable", "rate"]}]'
$ws.Range("G8").Value = '[{"factor": ["ct", "compute tomography"], "outcome": ["tbi", "patient", "death"]}, {"factor": ["marshall and rotterdam scoring system"], "outcome": ["tbi", "patient", "death"]}, {"factor": ["basal cistern absence"], "outcome": ["tbi", "patient", "death"]}, {"factor": ["positive midline shift"], "outcome": ["tbi", "patient", "death"]}, {"factor": ["hemorrhagic mass volume"], "outcome": ["tbi", "patient", "death"]}, {"factor": ["intraventricular", "subarachnoid hemorrhage"], "outcome": ["tbi", "patient", "death"]}]'
$ws.Range("G9").Value = '[{"factor": ["apache ii"], "outcome": ["tbi", "icu-treated", "mortality", "patient", "six-month"]}, {"factor": ["sap ii"], "outcome": ["tbi", "icu-treated", "mortality", "patient", "six-month"]}, {"factor": ["sofa"], "outcome": ["tbi", "icu-treated", "mortality", "patient", "six-month"]}, {"factor": ["age"], "outcome": ["tbi", "icu-treated", "mortality", "patient", "six-month"]}, {"factor": ["glasgow coma scale"], "outcome": ["tbi", "icu-treated", "mortality", "patient", "six-month"]}]'
$ws.Range("G10").Value = '[{"factor": ["v/c ratio"], "outcome": ["lcf score", "gos score", "drs"]}]'
$ws.Range("G11").Value = '[{"factor": ["serum", "timp-1 level"], "outcome": ["tbi", "patient", "severe", "mortality"]}]'
$ws.Range("G12").Value = '[{"factor": ["evidence", "mri", "contusion"], "outcome": ["glasgow outcome scale-extended", "gos-e"]}, {"factor": ["fa", "roi", "reduce", "severely"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["age"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}, {"factor": ["fa", "roi", "reduce", "severely"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}]'
$ws.Range("G13").Value = '[{"factor": ["depressive symptom", "preinjury"], "outcome": ["health-related quality-of-life", "physical problem", "cognitive", "affective/behavioral"]}]'
$ws.Range("G15").Value = '[{"factor": ["age"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["sex"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["deyo-charlson comorbidity index", "score"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["hospital volume"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["volume", "surgeon"], "outcome": ["in-hospital mortality rate"]}]'
$ws.Range("G17").Value = '[{"factor": ["amyloid-beta1-42 (abeta42)", "cerebrospinal fluid", "concentration", "csf"], "outcome": ["mortality"]}, {"factor": ["amyloid-beta1-42 (abeta42)", "concentration", "plasma"], "outcome": ["mortality"]}, {"factor": ["abeta42", "concentration", "change", "csf"], "outcome": ["neurological status"]}]'
$ws.Range("G18").Value = '[{"factor": ["plasminogen activator receptor", "urokinase", "supar", "soluble"], "outcome": ["traumatic brain injury", "tbi"]}, {"factor": ["score", "glasgow coma scale"], "outcome": ["tbi", "severity"]}, {"factor": ["d-dimer"], "outcome": ["tbi", "patient", "mortality"]}]'
$ws.Range("G19").Value = '[{"factor": ["motor score", "gcs"], "outcome": ["month", "mortality"]}, {"factor": ["age"], "outcome": ["month", "mortality"]}, {"factor": ["sex"], "outcome": ["month", "mortality"]}, {"factor": ["injury", "mechanism"], "outcome": ["month", "mortality"]}, {"factor": ["glasgow coma scale"], "outcome": ["month", "mortality"]}, {"factor": ["intubation"], "outcome": ["month", "mortality"]}, {"factor": ["pupil"], "outcome": ["month", "mortality"]}, {"factor": ["systolic blood pressure"], "outcome": ["month", "mortality"]}, {"factor": ["respiratory rate"], "outcome": ["month", "mortality"]}, {"factor": ["body temperature"], "outcome": ["month", "mortality"]}, {"factor": ["ph", "arterial"], "outcome": ["month", "mortality"]}, {"factor": ["arterial partial pressure", "carbon dioxide"], "outcome": ["month", "mortality"]}, {"factor": ["arterial partial pressure"], "outcome": ["month", "mortality"]}, {"factor": ["serum sodium"], "outcome": ["month", "mortality"]}, {"factor": ["serum potassium"], "outcome": ["month", "mortality"]}, {"factor": ["serum chloride"], "outcome": ["month", "mortality"]}, {"factor": ["serum calcium"], "outcome": ["month", "mortality"]}, {"factor": ["serum glucose"], "outcome": ["month", "mortality"]}, {"factor": ["urea nitrogen"], "outcome": ["month", "mortality"]}, {"factor": ["creatinine"], "outcome": ["month", "mortality"]}, {"factor": ["ratio", "international"], "outcome": ["month", "mortality"]}]'
$ws.Range("G20").Value = '[{"factor": ["aptt"], "outcome": ["deterioration", "surgery"]}, {"factor": ["fdp"], "outcome": ["deterioration", "surgery"]}, {"factor": ["d-dimer"], "outcome": ["deterioration", "surgery"]}]'
$ws.Range("G21").Value = '[{"factor": ["central conduction time", "cct"], "outcome": ["clinical outcome", "long-term"]}, {"factor": ["latency"], "outcome": ["clinical outcome", "long-term"]}]'
$ws.Range("G22").Value = '[{"factor": ["systolic blood pressure", "sbp"], "outcome": ["mortality"]}]'
$ws.Range("G24").Value = '[{"factor": ["tsp-1"], "outcome": ["1-week", "outcome", "unfavorable", "mortality"]}]'
$ws.Range("G25").Value = '[{"factor": ["plasma level", "brain-derived neurotrophic factor (bdnf)"], "outcome": ["tbi", "severe", "mortality", "intensive care unit", "patient", "icu"]}]'
$ws.Range("G26").Value = '[{"factor": ["crash-ct model"], "outcome": ["days", "death"]}, {"factor": ["age"], "outcome": ["older", "performance", "patient", "model"]}, {"factor": ["glasgow coma scale score"], "outcome": ["discrimination", "model"]}, {"factor": ["hosmer-lemeshow p value"], "outcome": ["calibration", "model"]}]'
$ws.Range("G27").Value = '[{"factor": ["time to death"], "outcome": ["withdrawal", "life-sustaining", "therapy"]}, {"factor": ["score", "glasgow coma scale"], "outcome": ["mortality"]}, {"factor": ["head abbreviate injury scale", "ais", "score"], "outcome": ["mortality"]}, {"factor": ["multiple", "comorbiditie"], "outcome": ["mortality"]}, {"factor": ["traumatic", "subarachnoid hemorrhage"], "outcome": ["mortality"]}, {"factor": ["intracerebral mass lesion"], "outcome": ["mortality"]}, {"factor": ["brainstem lesion"], "outcome": ["mortality"]}, {"factor": ["absent", "basal cistern", "sign of compress"], "outcome": ["mortality"]}]'
$ws.Range("G29").Value = '[{"factor": ["score", "glasgow coma scale"], "outcome": ["1-week", "outcome", "unfavorable", "mortality"]}]'
$ws.Range("G30").Value = '[{"factor": ["gcs", "glasgow coma scale"], "outcome": ["tbi", "severity"]}, {"factor": ["duration", "pta", "post-traumatic amnesia"], "outcome": ["tbi", "olfactory problem"]}]'
$ws.Range("G31").Value = '[{"factor": ["level", "il-6"], "outcome": ["development", "septic"]}, {"factor": ["c-reactive protein level"], "outcome": ["development", "multiple organ dysfunction"]}]'
$ws.Range("G32").Value = '[{"factor": ["rotterdam"], "outcome": ["hospital discharge", "death", "weeks"]}, {"factor": ["age"], "outcome": ["hospital discharge", "death", "weeks"]}, {"factor": ["sex"], "outcome": ["hospital discharge", "death", "weeks"]}, {"factor": ["glasgow coma scale score"], "outcome": ["hospital discharge", "death", "weeks"]}]'
$ws.Range("G33").Value = '[{"factor": ["glucose", "serum", "admission", "level of"], "outcome": ["traumatic brain injury", "patient", "severe", "outcome"]}]'
$ws.Range("G35").Value = '[{"factor": ["age"], "outcome": ["poor outcome"]}, {"factor": ["glasgow coma scale"], "outcome": ["poor outcome"]}, {"factor": ["severity score", "injury"], "outcome": ["poor outcome"]}, {"factor": ["ais", "head"], "outcome": ["poor outcome"]}]'
$ws.Range("G36").Value = '[{"factor": ["il-6"], "outcome": ["gos", "favorable", "year"]}, {"factor": ["pg"], "outcome": ["gos", "favorable", "year"]}, {"factor": ["gfap"], "outcome": ["gos", "unfavorable", "year", "score"]}, {"factor": ["pg"], "outcome": ["year", "survival status"]}, {"factor": ["gfap"], "outcome": ["year", "survival status"]}, {"factor": ["il-6"], "outcome": ["year", "survival status"]}]'
$ws.Range("G37").Value = '[{"factor": ["acute"], "outcome": ["score", "drs"]}, {"factor": ["subacute", "fa"], "outcome": ["score", "drs"]}]'
$ws.Range("G38").Value = '[{"factor": ["glasgow coma scale"], "outcome": ["mortality"]}, {"factor": ["mechanical ventilation"], "outcome": ["neurological"]}, {"factor": ["blood transfusion"], "outcome": ["neurological"]}, {"factor": ["neurosurgical intervention"], "outcome": ["neurological"]}, {"factor": ["concomitant", "injury"], "outcome": ["non-neurological", "complication"]}]'
$ws.Range("G39").Value = '[{"factor": ["traumatic brain injury"], "outcome": ["mortality"]}, {"factor": ["systolic blood pressure", "refer", "hospital"], "outcome": ["mortality"]}]'
